$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1470.6316
$ws.Range("I4").Value = 1071.7778
$ws.Range("K4").Value = 1071.7778
$ws.Range("M4").Value = -957.7778000000001
$ws.Range("H9").Value = 438.875
$ws.Range("I9").Value = 377
$ws.Range("K9").Value = 377
$ws.Range("M9").Value = -208
$ws.Range("H40").Value = 62501932
$ws.Range("J40").Value = 100001896
$ws.Range("L40").Value = 100001896
$ws.Range("N40").Value = -100002246
$ws.Range("H94").Value = 3324.1667
$ws.Range("I94").Value = 3471.818
$ws.Range("K94").Value = 3471.818
$ws.Range("M94").Value = -3020.818
$ws.Range("H97").Value = 1098.6666
$ws.Range("J97").Value = 1173.5
$ws.Range("L97").Value = 3520.5
$ws.Range("N97").Value = -4512.5
$ws.Range("H101").Value = 710.5833
$ws.Range("I101").Value = 235.5
$ws.Range("K101").Value = 706.5
$ws.Range("M101").Value = 915.5
$ws.Range("H125").Value = 1260.6666
$ws.Range("J125").Value = 1521.3636
$ws.Range("L125").Value = 13692.2724
$ws.Range("N125").Value = -18612.2724
$ws.Range("H132").Value = 1926.262
$ws.Range("I132").Value = 1997.079
$ws.Range("K132").Value = 5991.237
$ws.Range("M132").Value = -3461.237
$ws.Range("H138").Value = 4459.3945
$ws.Range("I138").Value = 2386.158
$ws.Range("J138").Value = 6532.6313
$ws.Range("K138").Value = 7158.474
$ws.Range("L138").Value = 19597.8939
$ws.Range("M138").Value = -2018.474
$ws.Range("N138").Value = -29877.8939

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10103.339
$ws.Range("I32").Value = 9737.804
$ws.Range("K32").Value = 9737.804
$ws.Range("M32").Value = -9450.804
$ws.Range("H60").Value = 200073660
$ws.Range("J60").Value = 500050000
$ws.Range("L60").Value = 500050000
$ws.Range("N60").Value = -500051466
$ws.Range("H61").Value = 21521796
$ws.Range("I61").Value = 25006620
$ws.Range("K61").Value = 25006620
$ws.Range("M61").Value = -25006408
$ws.Range("H122").Value = 2921.7878
$ws.Range("I122").Value = 2504.8147
$ws.Range("K122").Value = 7514.4441
$ws.Range("M122").Value = -5064.4441
$ws.Range("H136").Value = 21521796
$ws.Range("I136").Value = 25006620
$ws.Range("K136").Value = 75019860
$ws.Range("M136").Value = -75017310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2455.7778
$ws.Range("I94").Value = 3174.818
$ws.Range("J94").Value = 1325.8572
$ws.Range("K94").Value = 3174.818
$ws.Range("L94").Value = 1325.8572
$ws.Range("M94").Value = -2723.818
$ws.Range("N94").Value = -2227.8572
$ws.Range("H107").Value = 3439.4783
$ws.Range("I107").Value = 3797.2778
$ws.Range("J107").Value = 2151.4
$ws.Range("K107").Value = 3797.2778
$ws.Range("L107").Value = 2151.4
$ws.Range("M107").Value = -1877.2778
$ws.Range("N107").Value = -5991.4
$ws.Range("H134").Value = 4001369.5
$ws.Range("I134").Value = 1426.6666
$ws.Range("K134").Value = 4279.9998
$ws.Range("M134").Value = -1744.9998
$ws.Range("H139").Value = 116635.37
$ws.Range("J139").Value = 116635.37
$ws.Range("L139").Value = 116635.37
$ws.Range("N139").Value = -126915.37

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1765
$ws.Range("I58").Value = 1235.7333
$ws.Range("K58").Value = 1235.7333
$ws.Range("M58").Value = -1032.7333
$ws.Range("H103").Value = 27196.637
$ws.Range("J103").Value = 44397.4
$ws.Range("L103").Value = 44397.4
$ws.Range("N103").Value = -46741.4
$ws.Range("H122").Value = 3962.2632
$ws.Range("J122").Value = 5852.1113
$ws.Range("L122").Value = 17556.3339
$ws.Range("N122").Value = -22456.3339
$ws.Range("H134").Value = 1847.7222
$ws.Range("I134").Value = 1647.1562
$ws.Range("K134").Value = 4941.4686
$ws.Range("M134").Value = -2406.4686
$ws.Range("H136").Value = 1765
$ws.Range("I136").Value = 1235.7333
$ws.Range("K136").Value = 3707.199900000001
$ws.Range("M136").Value = -1157.199900000001
$ws.Range("H137").Value = 90000
$ws.Range("J137").Value = 90000
$ws.Range("L137").Value = 90000
$ws.Range("N137").Value = -100200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2049.889
$ws.Range("I68").Value = 738.5
$ws.Range("J68").Value = 3099
$ws.Range("K68").Value = 2215.5
$ws.Range("L68").Value = 9297
$ws.Range("M68").Value = -1404.5
$ws.Range("N68").Value = -10919
$ws.Range("H71").Value = 2049.889
$ws.Range("I71").Value = 738.5
$ws.Range("J71").Value = 3099
$ws.Range("K71").Value = 6646.5
$ws.Range("L71").Value = 27891
$ws.Range("M71").Value = -2590.5
$ws.Range("N71").Value = -36003
$ws.Range("H92").Value = 189.4
$ws.Range("J92").Value = 174.25
$ws.Range("L92").Value = 522.75
$ws.Range("N92").Value = -3018.75
$ws.Range("H132").Value = 3573.8286
$ws.Range("I132").Value = 2742.5386
$ws.Range("K132").Value = 24682.8474
$ws.Range("M132").Value = -22152.8474
$ws.Range("H134").Value = 14293396
$ws.Range("I134").Value = 17650666
$ws.Range("K134").Value = 52951998
$ws.Range("M134").Value = -52946928

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3239.5
$ws.Range("I97").Value = 903.1429000000001
$ws.Range("J97").Value = 8691
$ws.Range("K97").Value = 903.1429000000001
$ws.Range("L97").Value = 8691
$ws.Range("M97").Value = -407.1429000000001
$ws.Range("N97").Value = -9683
$ws.Range("H102").Value = 1150.1
$ws.Range("I102").Value = 1056.2593
$ws.Range("K102").Value = 1056.2593
$ws.Range("M102").Value = 565.7407000000001
$ws.Range("H132").Value = 6163196
$ws.Range("I132").Value = 5103.32
$ws.Range("K132").Value = 15309.96
$ws.Range("M132").Value = -12779.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1197.9412
$ws.Range("J46").Value = 1383.4445
$ws.Range("L46").Value = 1383.4445
$ws.Range("N46").Value = -1759.4445
$ws.Range("H132").Value = 3632.6
$ws.Range("J132").Value = 7288.4
$ws.Range("L132").Value = 21865.2
$ws.Range("N132").Value = -26925.2
$ws.Range("H136").Value = 3808.3635
$ws.Range("I136").Value = 3589.2
$ws.Range("K136").Value = 10767.6
$ws.Range("M136").Value = -8217.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 28006
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents() | Out-Null
$ws.Range("H70").Value = 39663.332
$ws.Range("J70").Value = 46995
$ws.Range("L70").Value = 46995
$ws.Range("N70").Value = -47625
$ws.Range("H73").Value = 39663.332
$ws.Range("J73").Value = 46995
$ws.Range("L73").Value = 46995
$ws.Range("N73").Value = -49179
$ws.Range("H113").Value = 818.7
$ws.Range("I113").Value = 603.7368
$ws.Range("J113").Value = 1190
$ws.Range("K113").Value = 1811.2104
$ws.Range("L113").Value = 3570
$ws.Range("M113").Value = 358.7896000000001
$ws.Range("N113").Value = -7910
$ws.Range("H132").Value = 835045.5600000001
$ws.Range("I132").Value = 1867.909
$ws.Range("J132").Value = 10000000
$ws.Range("K132").Value = 5603.727000000001
$ws.Range("L132").Value = 30000000
$ws.Range("M132").Value = -3073.727000000001
$ws.Range("N132").Value = -30005060
$ws.Range("H135").Value = 93540.375
$ws.Range("J135").Value = 93540.375
$ws.Range("L135").Value = 93540.375
$ws.Range("N135").Value = -103680.375
